$wb = $excel.ActiveWorkbook

# --- "Steps" sheet: insert a new "ALEX" row (Update_Nightly workflow)
#     right before the existing "Update Customer Unit Output" row ---
$wsSteps = $wb.Worksheets.Item("Steps")
$wsSteps.Range("A6:N6").Insert()
$wsSteps.Range("A6").Value = "Update_Nightly"
$wsSteps.Range("B6").Value = "ALEX"
$wsSteps.Range("C6").Value = "Formula"
$wsSteps.Range("D6").Value = "UpdateWorkflow (LIB_EWS)"
$wsSteps.Range("A6:D6").Style = "Normal"

# --- "Parameters" sheet: insert the matching ALEX parameter row ---
$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("A6:E6").Insert()
$wsParams.Range("A6").Value = "Update_Nightly"
$wsParams.Range("B6").Value = "ALEX"
$wsParams.Range("C6").Value = "Param1"
$wsParams.Range("D6").Value = "COUNTERPARTY_ALEX"
$wsParams.Range("A6:D6").Style = "Normal"

# --- Update each sheet's remembered selection, then leave
#     "Global Variables" as the active tab ---
$wsWorkflow = $wb.Worksheets.Item("Workflow")
$wsWorkflow.Range("A10").Select() | Out-Null

$wsSteps.Range("D6").Select() | Out-Null

$wsParams.Range("D10").Select() | Out-Null

$wsGlobalVars = $wb.Worksheets.Item("Global Variables")
$wsGlobalVars.Range("B11").Select() | Out-Null
$wsGlobalVars.Activate() | Out-Null
